$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows(310).Insert()
$ws.Rows(14).Copy()
$ws.Rows(310).PasteSpecial(-4122)
Write-Host "done"
